# Applies the "fixed a few things in example metadata files, added a
# bacteria fasta for Cd and then also changed a few of the test params"
# commit to the Cdiphtheriae_test_1 metadata sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: update the existing CP040557 record -----------------------
# id columns: sample_name / ncbi-spuid / isolate all carried the old
# "CP040557.1" value -> now just "CP040557"
$ws.Range("A3").Value = "CP040557"
$ws.Range("B3").Value = "CP040557"
$ws.Range("K3").Value = "CP040557"

# author changed from the placeholder "Michael" to "John Doe"
$ws.Range("E3").Value = "John Doe"

# newly populated submitting-lab columns
$ws.Range("F3").Value = "Black Bird Labs"
$ws.Range("G3").Value = "Bio intelligence"

# isolation_source was blank before, now "Clinical"
$ws.Range("L3").Value = "Clinical"

# lat_lon corrected from "NA" to "Not provided"
$ws.Range("V3").Value = "Not provided"

# new fasta_file_name-style + file_location + illumina fastq paths
$ws.Range("AE3").Value = "CP040557_Cd"
$ws.Range("AL3").Value = "local"
$ws.Range("AM3").Value = "/scicomp/instruments-pure/23-4-631_Illumina-MiSeq-M03083/220622_M03083_0094_000000000-KG73J/Alignment_1/20220624_225908/Fastq/CP040557_test_R1.fastq.gz"
$ws.Range("AN3").Value = "/scicomp/instruments-pure/23-4-631_Illumina-MiSeq-M03083/220622_M03083_0094_000000000-KG73J/Alignment_1/20220624_225908/Fastq/CP040557_test_R2.fastq.gz"

# --- Row 4: brand-new BX248355 record, mirrors row 3 -------------------
$ws.Range("A4").Value = "BX248355"
$ws.Range("B4").Value = "BX248355"
$ws.Range("K4").Value = "BX248355"

$ws.Range("E4").Value = "John Doe"
$ws.Range("F4").Value = "Black Bird Labs"
$ws.Range("G4").Value = "Bio intelligence"

$ws.Range("L4").Value = "Clinical"

$ws.Range("O4").Value = "Homo sapiens"
$ws.Range("P4").Value = "Cdiphtheriae"
$ws.Range("Q4").Value = "2022-06"
$ws.Range("R4").Value = "USA"
$ws.Range("T4").Value = "Not provided"
$ws.Range("V4").Value = "Not provided"

$ws.Range("AE4").Value = "BX248355_Cd"
$ws.Range("AL4").Value = "local"
$ws.Range("AM4").Value = "/scicomp/instruments-pure/23-4-631_Illumina-MiSeq-M03083/220622_M03083_0094_000000000-KG73J/Alignment_1/20220624_225908/Fastq/BX248355_test_R1.fastq.gz"
$ws.Range("AN4").Value = "/scicomp/instruments-pure/23-4-631_Illumina-MiSeq-M03083/220622_M03083_0094_000000000-KG73J/Alignment_1/20220624_225908/Fastq/BX248355_test_R2.fastq.gz"

# --- View/zoom cosmetics (best-effort) ---------------------------------
$excel.ActiveWindow.Zoom = 290
$ws.Range("AN5").Select()
